# Update the benchmark results table. The table has a single column;
# each row holds one metric value. Rows are addressed by (1-based) index
# so that duplicate values elsewhere in the table are not disturbed.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "3600"
$t.Cell(5, 1).Range.Text  = "0.00001"
$t.Cell(6, 1).Range.Text  = "0.00260"
$t.Cell(7, 1).Range.Text  = "0.00018"
$t.Cell(8, 1).Range.Text  = "0.00006"
$t.Cell(9, 1).Range.Text  = "0.00033"
$t.Cell(10, 1).Range.Text = "0.00039"
$t.Cell(11, 1).Range.Text = "0.00055"
$t.Cell(12, 1).Range.Text = "0.80313"

$t.Cell(44, 1).Range.Text = "99.84"
$t.Cell(45, 1).Range.Text = "0.8"
$t.Cell(46, 1).Range.Text = "495"
